$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.881.60'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.637.35'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.06'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.11'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.627'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.45%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.634.62'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.118'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.55%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.79'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.383'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.46%  '
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.61'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.116.04'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000184'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.794.58'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.651.11'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.59%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.14'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.80'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.52'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '345.05'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.71%  '
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.13'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('E25').Value = '  +5.84%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000111'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.33'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.64'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '582.16'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.18'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.161'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.23%  '
$ws.Range('E32').Value = '  -0.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.06'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.73'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.18%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.60'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.47'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.402'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.75'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.51%  '
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.91'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '152.83'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.74%  '
$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.54'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '163.14'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '24.28'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.90'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0589'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.08%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.634'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('E50').Value = '  -3.09%  '
$ws.Range('E51').Value = '  -2.02%  '
